# Refresh the cryptocurrency price/volume snapshot (GitHub Actions data pull).
# For each touched cell we set the new text value. A handful of new values look
# numeric to Excel (e.g. "540.96", "1.00") and would silently be coerced to a
# number on assignment, so those are written with a leading apostrophe (forces
# text entry, exactly like a user typing '540.96 into the cell) and then have
# their Style reset to "Normal" so they end up identical to the untouched cells
# (default style, plain text) instead of picking up a quote-prefix style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.070.14"
$ws.Range("E2").Value = "  -3.39%  "
# Row 3
$ws.Range("D3").Value = "3.284.50"
$ws.Range("E3").Value = "  -5.58%  "
# Row 4
$ws.Range("E4").Value = "  -0.05%  "
# Row 5
$ws.Range("D5").Value = "'540.96"
$ws.Range("E5").Value = "  -2.21%  "
# Row 6
$ws.Range("D6").Value = "'169.65"
$ws.Range("E6").Value = "  -4.71%  "
# Row 7
$ws.Range("E7").Value = "  -4.62%  "
# Row 8
$ws.Range("E8").Value = "  +0.04%  "
# Row 9
$ws.Range("D9").Value = "3.274.97"
$ws.Range("E9").Value = "  -5.70%  "
# Row 10
$ws.Range("D10").Value = "'0.604"
$ws.Range("E10").Value = "  -4.28%  "
# Row 11
$ws.Range("E11").Value = "  -1.18%  "
# Row 12
$ws.Range("D12").Value = "'52.27"
$ws.Range("E12").Value = "  -2.30%  "
# Row 13
$ws.Range("E13").Value = "  -2.92%  "
# Row 14
$ws.Range("D14").Value = "'8.76"
$ws.Range("E14").Value = "  -4.95%  "
# Row 15
$ws.Range("D15").Value = "3.819.27"
$ws.Range("E15").Value = "  -5.37%  "
# Row 16
$ws.Range("D16").Value = "'17.77"
$ws.Range("E16").Value = "  -4.87%  "
# Row 17
$ws.Range("E17").Value = "  -4.44%  "
# Row 18
$ws.Range("D18").Value = "3.288.26"
$ws.Range("E18").Value = "  -5.69%  "
# Row 19
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "63.070.71"
$ws.Range("E19").Value = "  -3.67%  "
# Row 20
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'11.52"
$ws.Range("E20").Value = "  -4.62%  "
# Row 21
$ws.Range("D21").Value = "'0.962"
$ws.Range("E21").Value = "  -2.57%  "
# Row 22
$ws.Range("D22").Value = "'410.79"
$ws.Range("E22").Value = "  -1.34%  "
# Row 23
$ws.Range("E23").Value = "  +6.02%  "
# Row 24
$ws.Range("D24").Value = "'3.99"
$ws.Range("E24").Value = "  -0.93%  "
# Row 25
$ws.Range("D25").Value = "'13.49"
$ws.Range("E25").Value = "  +5.34%  "
# Row 26
$ws.Range("D26").Value = "'82.25"
$ws.Range("E26").Value = "  -4.40%  "
# Row 27
$ws.Range("D27").Value = "'10.42"
$ws.Range("E27").Value = "  -3.17%  "
# Row 28
$ws.Range("E28").Value = "  -5.18%  "
# Row 29
$ws.Range("D29").Value = "'8.50"
$ws.Range("E29").Value = "  -5.77%  "
# Row 30
$ws.Range("D30").Value = "'28.70"
$ws.Range("E30").Value = "  -4.90%  "
# Row 31
$ws.Range("D31").Value = "'6.27"
$ws.Range("E31").Value = "  -3.33%  "
# Row 32
$ws.Range("D32").Value = "'11.22"
$ws.Range("E32").Value = "  -4.20%  "
# Row 33
$ws.Range("D33").Value = "'567.96"
$ws.Range("E33").Value = "  -6.86%  "
# Row 34
$ws.Range("E34").Value = "  -3.97%  "
# Row 35
$ws.Range("D35").Value = "'57.40"
$ws.Range("E35").Value = "  -3.25%  "
# Row 36
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.23%  "
# Row 37
$ws.Range("E37").Value = "  -1.11%  "
# Row 38
$ws.Range("D38").Value = "'34.74"
$ws.Range("E38").Value = "  -6.86%  "
# Row 39
$ws.Range("D39").Value = "'3.37"
$ws.Range("E39").Value = "  +4.52%  "
# Row 40
$ws.Range("D40").Value = "0.0₃0726"
$ws.Range("E40").Value = "  -7.14%  "
# Row 41
$ws.Range("D41").Value = "'0.360"
$ws.Range("E41").Value = "  -4.93%  "
# Row 42
$ws.Range("D42").Value = "3.095.74"
$ws.Range("E42").Value = "  -7.97%  "
# Row 43
$ws.Range("E43").Value = "  -0.30%  "
# Row 44
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").Value = "'3.22"
$ws.Range("E44").Value = "  +0.78%  "
# Row 45
$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").Value = "'2.73"
$ws.Range("E45").Value = "  -3.24%  "
# Row 46
$ws.Range("E46").Value = "  -4.02%  "
# Row 47
$ws.Range("D47").Value = "'2.38"
$ws.Range("E47").Value = "  -5.72%  "
# Row 48
$ws.Range("E48").Value = "  -3.99%  "
# Row 49
$ws.Range("E49").Value = "  -4.01%  "
# Row 50
$ws.Range("D50").Value = "'131.99"
$ws.Range("E50").Value = "  -4.06%  "
# Row 51
$ws.Range("D51").Value = "'7.93"
$ws.Range("E51").Value = "  -6.45%  "

# Cells written with a quote-prefix above need their style reset to "Normal" so
# they stay on the workbook default style (matches every other data cell) instead
# of acquiring a new quote-prefixed cell style.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"

